# Updates cryptocurrency price/volume data in cryptos.xlsx (Sheet1),
# plus a row swap (rows 48/49: THORChain <-> ordi) to reflect a new ranking snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are stored as plain text (e.g. "43.685.89" using
# "." as both a thousands- and decimal-like separator). Assigning a plain numeric-
# looking string via .Value would make Excel coerce it to a real Number, which
# would change the cell's stored type. Forcing the NumberFormat to Text ("@") before
# the write keeps it a string; resetting the range style to "Normal" afterwards
# avoids leaving a stray text format applied to the cells.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.685.89"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "2.247.25"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "323.37"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "101.44"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "0.580"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.556"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "37.29"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "0.0833"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").Value = "7.71"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("D14").Value = "2.589.43"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "0.859"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "14.25"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "2.245.86"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "43.628.21"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "13.72"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").Value = "0.0₃0987"
$ws.Range("E20").Value = "  +2.89%  "
$ws.Range("D21").Value = "6.58"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "65.27"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "236.80"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "10.14"
$ws.Range("E27").Value = "  +1.52%  "
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Value = "37.04"
$ws.Range("E29").Value = "  +7.65%  "
$ws.Range("D30").Value = "6.31"
$ws.Range("E30").Value = "  -1.36%  "
$ws.Range("D31").Value = "160.65"
$ws.Range("E31").Value = "  +4.17%  "
$ws.Range("D32").Value = "20.20"
$ws.Range("E32").Value = "  -1.28%  "
$ws.Range("D33").Value = "0.0855"
$ws.Range("E33").Value = "  -1.77%  "
$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").Value = "3.16"
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +9.14%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").Value = "3.78"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").Value = "4.27"
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "15.79"
$ws.Range("E41").Value = "  +20.94%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").Value = "1.806.82"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "0.200"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").Value = "82.28"
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("D47").Value = "1.71"
$ws.Range("E47").Value = "  +6.01%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "74.63"
$ws.Range("E48").Value = "  -2.90%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "5.21"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "58.73"
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "103.44"
$ws.Range("E51").Value = "  +0.28%  "

$priceRange.Style = "Normal"
